$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Portfolio_3"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2024-09-16"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = '[{''id'': 32981, ''symbol'': ''CATE'', ''name'': ''Cate'', ''price'': 0.029139323930875163, ''change_24h'': -36.08711097, ''change_7d'': 294.18440566, ''market_cap'': 29139323.930875164}, {''id'': 32902, ''symbol'': ''BRUH'', ''name'': ''Bruh'', ''price'': 0.0002525611914506638, ''change_24h'': -20.74363296, ''change_7d'': -64.96657158, ''market_cap'': 252561.1914506638}, {''id'': 32862, ''symbol'': ''MUNCAT'', ''name'': ''MUNCAT'', ''price'': 0.019374330182423315, ''change_24h'': -4.60176334, ''change_7d'': 8.75089258, ''market_cap'': 19374330.182423316}, {''id'': 32844, ''symbol'': ''WIWI'', ''name'': ''Wiggly Willy'', ''price'': 1.1970992615880443e-06, ''change_24h'': -6.43596478, ''change_7d'': -30.49371301, ''market_cap'': 779073.1964659699}, {''id'': 32698, ''symbol'': ''DOGS'', ''name'': ''DOGS'', ''price'': 0.000974934022244226, ''change_24h'': -8.84416531, ''change_7d'': -4.42049895, ''market_cap'': 503797155.99470377}, {''id'': 32821, ''symbol'': ''BULLS'', ''name'': ''BULLS'', ''price'': 0.0003321096010086998, ''change_24h'': 10.31073071, ''change_7d'': -69.08446803, ''market_cap'': 332109.6010086998}, {''id'': 32797, ''symbol'': ''IVfun'', ''name'': ''Invest Zone'', ''price'': 0.01509359367209487, ''change_24h'': -17.72769307, ''change_7d'': -43.50275799, ''market_cap'': 14191282.129985562}, {''id'': 32790, ''symbol'': ''EBULL'', ''name'': ''ETHEREUM IS GOOD'', ''price'': 0.0002195472998804108, ''change_24h'': -11.28277107, ''change_7d'': -47.72450816, ''market_cap'': 2195472.998804108}, {''id'': 32732, ''symbol'': ''SUNPEPE'', ''name'': ''sunpepe'', ''price'': 0.0007073959010905591, ''change_24h'': -22.32759297, ''change_7d'': -15.13403981, ''market_cap'': 707395.9010905592}, {''id'': 32727, ''symbol'': ''FOFAR'', ''name'': ''Fofar (fofar.meme)'', ''price'': 0.004762916399806206, ''change_24h'': -20.60774427, ''change_7d'': -22.81827769, ''market_cap'': 4524734.2197121}, {''id'': 32724, ''symbol'': ''CAT'', ''name'': "Simon''s Cat", ''price'': 3.128230254831375e-05, ''change_24h'': -10.57951281, ''change_7d'': 12.07102345, ''market_cap'': 211154144.82625225}, {''id'': 32717, ''symbol'': ''SUNDOG'', ''name'': ''SUNDOG'', ''price'': 0.32350882367672973, ''change_24h'': -6.11142331, ''change_7d'': 11.17595893, ''market_cap'': 322697710.70768976}, {''id'': 32638, ''symbol'': ''CATDOG'', ''name'': ''CatDog'', ''price'': 6.376179414161324e-05, ''change_24h'': -13.70645916, ''change_7d'': 10.80027526, ''market_cap'': 6376179.4141613245}, {''id'': 32625, ''symbol'': ''FUKU'', ''name'': ''Fuku-Kun'', ''price'': 1.8667628802515917e-06, ''change_24h'': 17.41560241, ''change_7d'': -17.07599284, ''market_cap'': 600761.2023986556}, {''id'': 32618, ''symbol'': ''WDOG'', ''name'': ''Wrapped Dog'', ''price'': 0.013036393535211509, ''change_24h'': -1.46995629, ''change_7d'': -11.93829317, ''market_cap'': 13036373.185401201}, {''id'': 32598, ''symbol'': ''GODCAT'', ''name'': ''Godcat Exploding Kittens'', ''price'': 0.00012854647107292624, ''change_24h'': -10.03226685, ''change_7d'': 1.96281608, ''market_cap'': 102442.33827401283}, {''id'': 32524, ''symbol'': ''GINNAN'', ''name'': ''Ginnan The Cat'', ''price'': 3.1159333185215106e-06, ''change_24h'': 6.07614516, ''change_7d'': -12.29124678, ''market_cap'': 21499934.79448856}, {''id'': 32521, ''symbol'': ''NEIRO'', ''name'': ''First Neiro On Ethereum'', ''price'': 0.00031501157832759856, ''change_24h'': 650.9717983, ''change_7d'': 573.24760554, ''market_cap'': 131812785.08008216}, {''id'': 32495, ''symbol'': ''$YAWN'', ''name'': ''YAWN'', ''price'': 3.590853100253668e-05, ''change_24h'': -2.61457542, ''change_7d'': -17.57092872, ''market_cap'': 230372.5100883174}, {''id'': 32464, ''symbol'': ''NEIRO'', ''name'': ''Neiro (neiro.lol)'', ''price'': 0.004294003486495922, ''change_24h'': 28.80036116, ''change_7d'': 10.66580436, ''market_cap'': 4294003.3533818135}, {''id'': 32462, ''symbol'': ''YOUNES'', ''name'': ''YOUNES'', ''price'': 0.0015603198612713744, ''change_24h'': -18.73639977, ''change_7d'': -5.37927022, ''market_cap'': 1560317.054255944}, {''id'': 32461, ''symbol'': ''NEIRO'', ''name'': ''Neiro Ethereum'', ''price'': 0.08901291215736677, ''change_24h'': -33.68366943, ''change_7d'': -46.54167583, ''market_cap'': 89012912.15736677}, {''id'': 32448, ''symbol'': ''BRAINLET'', ''name'': ''Brainlet'', ''price'': 0.013371055846627362, ''change_24h'': -12.36710061, ''change_7d'': -36.54155982, ''market_cap'': 12741353.20527166}, {''id'': 32415, ''symbol'': ''SKBDI'', ''name'': ''Skibidi Toilet'', ''price'': 0.04473303942948299, ''change_24h'': -10.90531428, ''change_7d'': -22.92623991, ''market_cap'': 3133042.6761646266}, {''id'': 32350, ''symbol'': ''SOY'', ''name'': ''Soyjak'', ''price'': 0.0006795624706939648, ''change_24h'': -2.64842734, ''change_7d'': -35.54658999, ''market_cap'': 629654.0362460715}, {''id'': 32310, ''symbol'': ''SAITAMA'', ''name'': ''SAITAMA INU'', ''price'': 0.00027430996870847366, ''change_24h'': -18.53043543, ''change_7d'': 23.72363344, ''market_cap'': 182820.29420731042}, {''id'': 32297, ''symbol'': ''BILL'', ''name'': ''BILL THE BEAR'', ''price'': 1.979811310241725e-06, ''change_24h'': -6.00642017, ''change_7d'': 5.29085747, ''market_cap'': 1759826.1496376554}, {''id'': 32288, ''symbol'': ''torsy'', ''name'': ''TORSY'', ''price'': 0.003005750585596339, ''change_24h'': -1.82467013, ''change_7d'': 6.5042258, ''market_cap'': 2104024.427036996}, {''id'': 32259, ''symbol'': ''MOBY'', ''name'': ''Moby'', ''price'': 0.00013450244292354262, ''change_24h'': -5.03668237, ''change_7d'': 11.54955673, ''market_cap'': 1345024.4292354262}, {''id'': 32128, ''symbol'': ''MSI'', ''name'': ''monkey shit inu'', ''price'': 1.8098729385003677e-06, ''change_24h'': 9.95461635, ''change_7d'': 45.91601572, ''market_cap'': 5328754.305610909}, {''id'': 32125, ''symbol'': ''COIN'', ''name'': ''8-Bit Coin'', ''price'': 0.0006319251620375989, ''change_24h'': -1.15136583, ''change_7d'': -15.6718967, ''market_cap'': 18957754.86112797}, {''id'': 32074, ''symbol'': ''FEG'', ''name'': ''FEED EVERY GORILLA'', ''price'': 0.00030113265003483, ''change_24h'': -8.35137175, ''change_7d'': -19.63470506, ''market_cap'': 24253684.03009346}, {''id'': 32019, ''symbol'': ''CRASH'', ''name'': ''Crash'', ''price'': 0.003306616540248982, ''change_24h'': -7.91456337, ''change_7d'': 3.14937666, ''market_cap'': 3234655.7588133174}, {''id'': 31923, ''symbol'': ''HAWKTUAH'', ''name'': ''Hawk Tuah'', ''price'': 0.0009020855931686253, ''change_24h'': 1.38498241, ''change_7d'': 11.43846528, ''market_cap'': 774638.6140107772}, {''id'': 31921, ''symbol'': ''RETARDIO'', ''name'': ''RETARDIO'', ''price'': 0.035894365968676796, ''change_24h'': -10.60680773, ''change_7d'': -5.89337012, ''market_cap'': 35021101.45103184}, {''id'': 31914, ''symbol'': ''BILLY'', ''name'': ''Billy'', ''price'': 0.02019622322252712, ''change_24h'': -0.77008517, ''change_7d'': -18.87591525, ''market_cap'': 18906445.087785527}, {''id'': 31908, ''symbol'': ''WAT'', ''name'': ''Wat'', ''price'': 6.434700983075277e-06, ''change_24h'': -18.66694023, ''change_7d'': -0.02011322, ''market_cap'': 2692184.956223106}, {''id'': 31847, ''symbol'': ''WOLF'', ''name'': ''LandWolf (SOL)'', ''price'': 0.0017975557046317865, ''change_24h'': -9.84772067, ''change_7d'': 13.26939103, ''market_cap'': 17975275.154204458}, {''id'': 31830, ''symbol'': ''DADDY'', ''name'': ''Daddy Tate'', ''price'': 0.07349639130174392, ''change_24h'': -2.40572335, ''change_7d'': 0.92716441, ''market_cap'': 44074928.089101754}, {''id'': 31798, ''symbol'': ''JENNER'', ''name'': ''Caitlyn Jenner (ETH)'', ''price'': 0.00019243234131918272, ''change_24h'': -12.99314125, ''change_7d'': -16.94863288, ''market_cap'': 185615.01420760492}, {''id'': 31770, ''symbol'': ''HONK'', ''name'': ''Pepoclown'', ''price'': 2.068912460203536e-09, ''change_24h'': -16.38580393, ''change_7d'': 49.56251288, ''market_cap'': 830159.8043793539}, {''id'': 31678, ''symbol'': ''MAGAA'', ''name'': ''MAGA AGAIN'', ''price'': 0.00687644569218029, ''change_24h'': -3.06701461, ''change_7d'': -10.03807586, ''market_cap'': 6876445.69218029}, {''id'': 31632, ''symbol'': ''PEIPEI'', ''name'': ''PeiPei (ETH)'', ''price'': 1.0483626872196531e-07, ''change_24h'': -5.82396976, ''change_7d'': -14.50373131, ''market_cap'': 43104096.930973984}, {''id'': 31569, ''symbol'': ''TROG'', ''name'': ''Trog'', ''price'': 1.9000862937687478e-06, ''change_24h'': -15.54342796, ''change_7d'': 24.62627876, ''market_cap'': 799347.3029255745}, {''id'': 31561, ''symbol'': ''SOLCAT'', ''name'': ''CatSolHat'', ''price'': 0.012440070426778086, ''change_24h'': -0.84177775, ''change_7d'': 6.55731762, ''market_cap'': 823022.171522676}, {''id'': 31510, ''symbol'': ''MOTHER'', ''name'': ''Mother Iggy'', ''price'': 0.05428805790894996, ''change_24h'': -8.61527654, ''change_7d'': 53.69338084, ''market_cap'': 53535796.66769585}, {''id'': 31496, ''symbol'': ''PAPU'', ''name'': ''Papu Token'', ''price'': 8.63254434919503e-09, ''change_24h'': 6.35814711, ''change_7d'': 85.70910873, ''market_cap'': 273697.5393746039}, {''id'': 31317, ''symbol'': ''MEOW'', ''name'': ''MeowCat'', ''price'': 0.007325005223657749, ''change_24h'': -8.74740746, ''change_7d'': -17.20302316, ''market_cap'': 732500.522365775}, {''id'': 31305, ''symbol'': ''MAGA'', ''name'': ''MAGA'', ''price'': 6.375211512978745e-05, ''change_24h'': -0.95305309, ''change_7d'': -27.25946893, ''market_cap'': 24879830.05801541}, {''id'': 31284, ''symbol'': ''HAMMY'', ''name'': ''SAD HAMSTER'', ''price'': 0.011266577684823262, ''change_24h'': -6.34657659, ''change_7d'': -4.08111829, ''market_cap'': 10612003.604557136}, {''id'': 31267, ''symbol'': ''UTYAB'', ''name'': ''UTYABSWAP'', ''price'': 0.0011104591165764781, ''change_24h'': -15.10572874, ''change_7d'': 51.251048, ''market_cap'': 1110459.1165764781}, {''id'': 31259, ''symbol'': ''FOXSY'', ''name'': ''Foxsy AI'', ''price'': 0.08529040891276468, ''change_24h'': -5.26974331, ''change_7d'': 7.83199334, ''market_cap'': 17334487.95459311}, {''id'': 31163, ''symbol'': ''SLOTH'', ''name'': ''Slothana'', ''price'': 0.00593063885816661, ''change_24h'': -1.49380992, ''change_7d'': 2.18500391, ''market_cap'': 10053133.581434773}, {''id'': 31152, ''symbol'': ''KENDU'', ''name'': ''Kendu Inu'', ''price'': 5.5953786970426215e-05, ''change_24h'': 2.67499453, ''change_7d'': 25.07916241, ''market_cap'': 53068146.50843015}, {''id'': 31121, ''symbol'': ''LABZ'', ''name'': ''Insane Labz'', ''price'': 0.017456755559443735, ''change_24h'': -9.01472971, ''change_7d'': 4.32482207, ''market_cap'': 864235.8817746732}, {''id'': 31106, ''symbol'': ''DOKY'', ''name'': ''Donkey King'', ''price'': 2.1838017714979863e-05, ''change_24h'': -1.07148149, ''change_7d'': -21.9366622, ''market_cap'': 189435.2379487126}, {''id'': 31061, ''symbol'': ''DUREV'', ''name'': ''Povel Durev'', ''price'': 0.03528661976631283, ''change_24h'': -2.0411243, ''change_7d'': -7.14347822, ''market_cap'': 3352228.877799719}, {''id'': 31051, ''symbol'': ''LOBO'', ''name'': ''LOBO•THE•WOLF•PUP'', ''price'': 0.000510869249058204, ''change_24h'': -6.46887891, ''change_7d'': -19.80979546, ''market_cap'': 8372013.705265315}, {''id'': 31044, ''symbol'': ''HEGE'', ''name'': ''Hege'', ''price'': 0.01528185216378374, ''change_24h'': -1.1678672, ''change_7d'': 3.00578079, ''market_cap'': 15279596.514572168}, {''id'': 31036, ''symbol'': ''CARLO'', ''name'': ''Carlo'', ''price'': 0.0008650162782807935, ''change_24h'': -11.76720204, ''change_7d'': -25.17629652, ''market_cap'': 801811.1096963773}, {''id'': 30979, ''symbol'': ''ELON'', ''name'': ''Elon MemeLord'', ''price'': 4.536133694681866e-06, ''change_24h'': -6.32402697, ''change_7d'': 12.44899525, ''market_cap'': 3108248.540217904}, {''id'': 30968, ''symbol'': ''LONG'', ''name'': ''Long'', ''price'': 1.0422856620233421e-06, ''change_24h'': -5.74407207, ''change_7d'': 23.08945372, ''market_cap'': 816853.693589553}, {''id'': 30953, ''symbol'': ''SHIBTC'', ''name'': ''ShibaBitcoin'', ''price'': 0.6970305061463759, ''change_24h'': -0.07714365, ''change_7d'': -0.15680246, ''market_cap'': 60266007.566240266}, {''id'': 30943, ''symbol'': ''$MICHI'', ''name'': ''michi (SOL)'', ''price'': 0.11288322873985437, ''change_24h'': -4.39640829, ''change_7d'': 2.13794535, ''market_cap'': 62737047.35219478}, {''id'': 30933, ''symbol'': ''DOG'', ''name'': ''Dog (Runes)'', ''price'': 0.0021400667517258573, ''change_24h'': -2.44059291, ''change_7d'': -14.80441103, ''market_cap'': 214006675.17258573}, {''id'': 30912, ''symbol'': ''MANEKI'', ''name'': ''MANEKI'', ''price'': 0.00565028992493853, ''change_24h'': 4.40875026, ''change_7d'': 30.09895904, ''market_cap'': 47591754.75285718}, {''id'': 30867, ''symbol'': ''WHY'', ''name'': ''WHY'', ''price'': 2.692710739893155e-07, ''change_24h'': 0.49866255, ''change_7d'': 4.01470162, ''market_cap'': 113093851.0755125}, {''id'': 30859, ''symbol'': ''HOPPY'', ''name'': ''Hoppy'', ''price'': 5.368310561811881e-05, ''change_24h'': -6.35613079, ''change_7d'': 1.05730674, ''market_cap'': 22583945.702486403}, {''id'': 30828, ''symbol'': ''CATA'', ''name'': ''CATAMOTO'', ''price'': 0.00026764432395804066, ''change_24h'': -3.47241182, ''change_7d'': 16.08540717, ''market_cap'': 5186322.50648218}, {''id'': 30647, ''symbol'': ''GOAT'', ''name'': ''Sonic The Goat'', ''price'': 3.1491388254117285e-07, ''change_24h'': -3.36795121, ''change_7d'': -6.36869642, ''market_cap'': 279046.08410720347}, {''id'': 30645, ''symbol'': ''POWSCHE'', ''name'': ''Powsche'', ''price'': 0.015752725261451144, ''change_24h'': -8.62687032, ''change_7d'': 21.47726974, ''market_cap'': 1497290.7706034856}, {''id'': 30629, ''symbol'': ''KITTENWIF'', ''name'': ''KittenWifHat'', ''price'': 0.0008044336871728509, ''change_24h'': -2.44670259, ''change_7d'': 2.28799903, ''market_cap'': 786397.7067117799}, {''id'': 30601, ''symbol'': ''MOUTAI'', ''name'': ''Moutai'', ''price'': 0.005497809688962206, ''change_24h'': -14.527483, ''change_7d'': -11.69149726, ''market_cap'': 3830419.7406890998}, {''id'': 30550, ''symbol'': ''HAMI'', ''name'': ''HAMI'', ''price'': 0.0016212390216391401, ''change_24h'': -7.37546275, ''change_7d'': 0.31642221, ''market_cap'': 1560015.4545347888}, {''id'': 30493, ''symbol'': ''NUB'', ''name'': ''nubcat'', ''price'': 0.008871408835112867, ''change_24h'': 3.39199983, ''change_7d'': -8.28930856, ''market_cap'': 8427840.637823658}, {''id'': 30484, ''symbol'': ''PUPS'', ''name'': ''PUPS (Ordinals)'', ''price'': 2.507201534089858, ''change_24h'': -5.80841652, ''change_7d'': -6.86848297, ''market_cap'': 19480955.919878196}, {''id'': 30409, ''symbol'': ''NSO'', ''name'': "NeverSurrenderOne''s", ''price'': 0.00018821080490818678, ''change_24h'': -2.97789603, ''change_7d'': 2.97455753, ''market_cap'': 121059.34096644583}, {''id'': 30407, ''symbol'': ''ROOST'', ''name'': ''Roost Coin'', ''price'': 0.0017491617219330038, ''change_24h'': -15.80491292, ''change_7d'': 24.08285034, ''market_cap'': 1639216.0121811486}, {''id'': 30402, ''symbol'': ''CAW'', ''name'': ''crow with knife'', ''price'': 3.6141720897967874e-08, ''change_24h'': -4.32496404, ''change_7d'': -0.92919839, ''market_cap'': 27824115.76693413}, {''id'': 30361, ''symbol'': ''DOGEMOB'', ''name'': ''DOGEMOB'', ''price'': 0.00034513339900209126, ''change_24h'': -0.77270971, ''change_7d'': -3.85510172, ''market_cap'': 172376.90788386718}, {''id'': 30309, ''symbol'': ''SC'', ''name'': ''Shark Cat'', ''price'': 0.011597929638090474, ''change_24h'': -3.60634632, ''change_7d'': -9.5873231, ''market_cap'': 11480738.578423053}, {''id'': 30303, ''symbol'': ''GB'', ''name'': ''Gary Banking'', ''price'': 9.603228198005154e-05, ''change_24h'': -0.71192099, ''change_7d'': -36.68830043, ''market_cap'': 1008338.9607905411}, {''id'': 30285, ''symbol'': ''MUMU'', ''name'': ''Mumu the Bull'', ''price'': 3.311633143858145e-05, ''change_24h'': -9.35557125, ''change_7d'': -5.18532352, ''market_cap'': 75611354.21490838}, {''id'': 30270, ''symbol'': ''FOMO'', ''name'': ''FOMO BULL CLUB'', ''price'': 9.32242505811339e-06, ''change_24h'': -1.72464896, ''change_7d'': -37.68931312, ''market_cap'': 311961.20007738535}, {''id'': 30218, ''symbol'': ''MABA'', ''name'': ''Make America Based Again'', ''price'': 0.00021710814260935348, ''change_24h'': -11.23953804, ''change_7d'': -5.23313165, ''market_cap'': 200261.91070827295}, {''id'': 30193, ''symbol'': ''BENJI'', ''name'': ''Basenji'', ''price'': 0.014691124333341283, ''change_24h'': -4.97266429, ''change_7d'': -19.38632383, ''market_cap'': 13545875.8260895}, {''id'': 30133, ''symbol'': ''COINYE'', ''name'': ''Coinye West'', ''price'': 0.0009798889975638088, ''change_24h'': 9.80172775, ''change_7d'': 31.7042877, ''market_cap'': 785130.8913634287}, {''id'': 30126, ''symbol'': ''MEW'', ''name'': ''cat in a dogs world'', ''price'': 0.0042210796171238745, ''change_24h'': -1.23964412, ''change_7d'': 1.715932, ''market_cap'': 375207077.0739257}, {''id'': 30119, ''symbol'': ''WIF'', ''name'': ''DOGWIFHOOD'', ''price'': 0.000458884941988761, ''change_24h'': 7.13712257, ''change_7d'': 11.69619647, ''market_cap'': 458389.4256385081}, {''id'': 30096, ''symbol'': ''DEGEN'', ''name'': ''Degen'', ''price'': 0.002952136717829938, ''change_24h'': -9.03085859, ''change_7d'': -1.365551, ''market_cap'': 41860144.016163304}, {''id'': 30071, ''symbol'': ''LADYF'', ''name'': ''Milady Wif Hat'', ''price'': 1.332583878802956e-05, ''change_24h'': -3.66171163, ''change_7d'': 2.52262203, ''market_cap'': 3331436.323696907}, {''id'': 30063, ''symbol'': ''GIGA'', ''name'': ''Gigachad'', ''price'': 0.0164742932798592, ''change_24h'': -1.80681799, ''change_7d'': 2.76960118, ''market_cap'': 153250661.65296072}, {''id'': 30011, ''symbol'': ''SKID'', ''name'': ''Success Kid'', ''price'': 0.025331225184864834, ''change_24h'': -4.81351327, ''change_7d'': 0.56263472, ''market_cap'': 2197298.8161807247}, {''id'': 30008, ''symbol'': ''APU'', ''name'': ''Apu Apustaja'', ''price'': 0.00018650718384005105, ''change_24h'': -0.72103611, ''change_7d'': 6.62478854, ''market_cap'': 59373911.31698716}, {''id'': 30007, ''symbol'': ''HUND'', ''name'': ''HUND'', ''price'': 0.006278814398861651, ''change_24h'': -5.22448462, ''change_7d'': 13.23043212, ''market_cap'': 2397421.5669479524}, {''id'': 29999, ''symbol'': ''CHKN'', ''name'': ''Chickencoin'', ''price'': 4.834561903759423e-08, ''change_24h'': -6.19775309, ''change_7d'': -8.30257855, ''market_cap'': 3105674.2213560157}, {''id'': 29932, ''symbol'': ''BSHIB'', ''name'': ''Based Shiba Inu'', ''price'': 2.0321708154334188e-05, ''change_24h'': -2.11104076, ''change_7d'': 6.98052083, ''market_cap'': 183736.5113856465}, {''id'': 29920, ''symbol'': ''SLERF'', ''name'': ''SLERF'', ''price'': 0.13511996743584379, ''change_24h'': -0.37994926, ''change_7d'': 3.28035753, ''market_cap'': 67559679.69799516}, {''id'': 29879, ''symbol'': ''ANDY'', ''name'': ''ANDY (ETH)'', ''price'': 6.502301923744789e-05, ''change_24h'': -10.48923845, ''change_7d'': 11.66785725, ''market_cap'': 65023019.23744789}, {''id'': 29870, ''symbol'': ''BOME'', ''name'': ''BOOK OF MEME'', ''price'': 0.005999715543528559, ''change_24h'': -4.49814322, ''change_7d'': -3.97460541, ''market_cap'': 413750959.9093467}]'
$ws.Range("D4").Value = 58743.11
$ws.Range("E4").Value = 132.09
$ws.Range("F4").Value = "[]"
$ws.Rows.Item(4).AutoFit()
